$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the mora period for the existing worker (RAFAEL DAVID ECHAVARRIA LLERENA) ---
$ws.Range("E16").Value = "2507"

# --- Insert two new rows for the new worker before the signature block ---
$ws.Range("A17:A18").EntireRow.Insert(-4121)

# Copy the formatting of the existing data row (16) into the two new rows
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J18").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# New worker, first overdue period
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "72346876"
$ws.Range("D17").Value = "JOHN EDINSON MEJIA SOLANO"
$ws.Range("E17").Value = "2212"
$ws.Range("F17").Value = 40000
$ws.Range("G17").Value = 1000000

# New worker, second overdue period
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "72346876"
$ws.Range("D18").Value = "JOHN EDINSON MEJIA SOLANO"
$ws.Range("E18").Value = "2211"
$ws.Range("F18").Value = 40000
$ws.Range("G18").Value = 1000000

# --- Update the summary block: totals now reflect 2 workers / 3 periods / new total due ---
$ws.Range("E11").Value = 136940
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 3
